$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 373.34375
$ws.Cells.Item(33, 9).Value = 325.44
$ws.Cells.Item(33, 10).Value = 544.4286
$ws.Cells.Item(33, 11).Value = 325.44
$ws.Cells.Item(33, 12).Value = 544.4286
$ws.Cells.Item(33, 13).Value = -96.44
$ws.Cells.Item(33, 14).Value = -1002.4286
$ws.Cells.Item(127, 8).Value = 2134.5454
$ws.Cells.Item(127, 9).Value = 1245
$ws.Cells.Item(127, 10).Value = 2642.8572
$ws.Cells.Item(127, 11).Value = 3735
$ws.Cells.Item(127, 12).Value = 7928.571599999999
$ws.Cells.Item(127, 13).Value = 1225
$ws.Cells.Item(127, 14).Value = -17848.5716
$ws.Cells.Item(131, 8).Value = 1000
$ws.Cells.Item(131, 10).Value = 1000
$ws.Cells.Item(131, 12).Value = 3000
$ws.Cells.Item(131, 14).Value = -13080
$ws.Cells.Item(135, 8).Value = 90911870
$ws.Cells.Item(135, 9).Value = 942
$ws.Cells.Item(135, 10).Value = 250006000
$ws.Cells.Item(135, 11).Value = 8478
$ws.Cells.Item(135, 12).Value = 2250054000
$ws.Cells.Item(135, 13).Value = -5943
$ws.Cells.Item(135, 14).Value = -2250059070
$ws.Cells.Item(137, 8).Value = 1330.4634
$ws.Cells.Item(137, 9).Value = 869.875
$ws.Cells.Item(137, 11).Value = 2609.625
$ws.Cells.Item(137, 13).Value = -59.625
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 1276.75
$ws.Cells.Item(122, 9).Value = 1012.25
$ws.Cells.Item(122, 10).Value = 2202.5
$ws.Cells.Item(122, 11).Value = 3036.75
$ws.Cells.Item(122, 12).Value = 6607.5
$ws.Cells.Item(122, 13).Value = -586.75
$ws.Cells.Item(122, 14).Value = -11507.5
$ws.Cells.Item(132, 8).Value = 2420.5881
$ws.Cells.Item(132, 9).Value = 1964.5834
$ws.Cells.Item(132, 10).Value = 3515
$ws.Cells.Item(132, 11).Value = 5893.7502
$ws.Cells.Item(132, 12).Value = 10545
$ws.Cells.Item(132, 13).Value = -3363.7502
$ws.Cells.Item(132, 14).Value = -15605
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 3045.3333
$ws.Cells.Item(134, 9).Value = 1457.3334
$ws.Cells.Item(134, 10).Value = 4633.3335
$ws.Cells.Item(134, 11).Value = 4372.0002
$ws.Cells.Item(134, 12).Value = 13900.0005
$ws.Cells.Item(134, 13).Value = -1837.0002
$ws.Cells.Item(134, 14).Value = -18970.0005
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2039.3125
$ws.Cells.Item(31, 9).Value = 2246.9092
$ws.Cells.Item(31, 10).Value = 1582.6
$ws.Cells.Item(31, 11).Value = 2246.9092
$ws.Cells.Item(31, 12).Value = 1582.6
$ws.Cells.Item(31, 13).Value = -1951.9092
$ws.Cells.Item(31, 14).Value = -2172.6
$ws.Cells.Item(34, 8).Value = 2039.3125
$ws.Cells.Item(34, 9).Value = 2246.9092
$ws.Cells.Item(34, 10).Value = 1582.6
$ws.Cells.Item(34, 11).Value = 2246.9092
$ws.Cells.Item(34, 12).Value = 1582.6
$ws.Cells.Item(34, 13).Value = -2044.9092
$ws.Cells.Item(34, 14).Value = -1986.6
$ws.Cells.Item(99, 8).Value = 1644.5555
$ws.Cells.Item(99, 9).Value = 1606.5
$ws.Cells.Item(99, 11).Value = 1606.5
$ws.Cells.Item(99, 13).Value = -108.5
$ws.Cells.Item(107, 8).Value = 521.8182
$ws.Cells.Item(107, 9).Value = 418.57144
$ws.Cells.Item(107, 11).Value = 418.57144
$ws.Cells.Item(107, 13).Value = 1501.42856
$ws.Cells.Item(126, 8).Value = 1644.5555
$ws.Cells.Item(126, 9).Value = 1606.5
$ws.Cells.Item(126, 11).Value = 4819.5
$ws.Cells.Item(126, 13).Value = -2349.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(41, 8).Value = 638.46155
$ws.Cells.Item(41, 10).Value = 638.46155
$ws.Cells.Item(41, 12).Value = 1915.38465
$ws.Cells.Item(41, 14).Value = -2591.38465
$ws.Cells.Item(64, 8).Value = 5000
$ws.Cells.Item(64, 9).Value = 0
$ws.Cells.Item(64, 10).Value = 5000
$ws.Cells.Item(64, 11).Value = 0
$ws.Cells.Item(64, 12).Value = 15000
$ws.Cells.Item(64, 13).ClearContents()
$ws.Cells.Item(64, 14).Value = -15540
$ws.Cells.Item(67, 8).Value = 5000
$ws.Cells.Item(67, 9).Value = 0
$ws.Cells.Item(67, 10).Value = 5000
$ws.Cells.Item(67, 11).Value = 0
$ws.Cells.Item(67, 12).Value = 15000
$ws.Cells.Item(67, 13).ClearContents()
$ws.Cells.Item(67, 14).Value = -16872
$ws.Cells.Item(86, 8).Value = 454.14285
$ws.Cells.Item(86, 9).Value = 296.66666
$ws.Cells.Item(86, 10).Value = 572.25
$ws.Cells.Item(86, 11).Value = 889.9999799999999
$ws.Cells.Item(86, 12).Value = 1716.75
$ws.Cells.Item(86, 13).Value = 296.0000200000001
$ws.Cells.Item(86, 14).Value = -4088.75
$ws.Cells.Item(89, 8).Value = 454.14285
$ws.Cells.Item(89, 9).Value = 296.66666
$ws.Cells.Item(89, 10).Value = 572.25
$ws.Cells.Item(89, 11).Value = 2669.99994
$ws.Cells.Item(89, 12).Value = 5150.25
$ws.Cells.Item(89, 13).Value = 3258.00006
$ws.Cells.Item(89, 14).Value = -17006.25
$ws.Cells.Item(92, 8).Value = 240.23077
$ws.Cells.Item(92, 9).Value = 234.60869
$ws.Cells.Item(92, 10).Value = 283.33334
$ws.Cells.Item(92, 11).Value = 703.82607
$ws.Cells.Item(92, 12).Value = 850.0000200000001
$ws.Cells.Item(92, 13).Value = 544.17393
$ws.Cells.Item(92, 14).Value = -3346.00002
$ws.Cells.Item(94, 8).Value = 4001.125
$ws.Cells.Item(94, 9).Value = 3612
$ws.Cells.Item(94, 10).Value = 4056.7144
$ws.Cells.Item(94, 11).Value = 10836
$ws.Cells.Item(94, 12).Value = 12170.1432
$ws.Cells.Item(94, 13).Value = -10160
$ws.Cells.Item(94, 14).Value = -13522.1432
$ws.Cells.Item(95, 8).Value = 6500
$ws.Cells.Item(95, 9).Value = 6500
$ws.Cells.Item(95, 10).Value = 6500
$ws.Cells.Item(95, 11).Value = 19500
$ws.Cells.Item(95, 12).Value = 19500
$ws.Cells.Item(95, 13).Value = -17441
$ws.Cells.Item(95, 14).Value = -23618
$ws.Cells.Item(96, 8).Value = 8966.666999999999
$ws.Cells.Item(96, 10).Value = 8966.666999999999
$ws.Cells.Item(96, 12).Value = 26900.001
$ws.Cells.Item(96, 14).Value = -31018.001
$ws.Cells.Item(97, 8).Value = 1083.3334
$ws.Cells.Item(97, 9).Value = 833.3333
$ws.Cells.Item(97, 10).Value = 1333.3334
$ws.Cells.Item(97, 11).Value = 2499.9999
$ws.Cells.Item(97, 12).Value = 4000.0002
$ws.Cells.Item(97, 13).Value = -2003.9999
$ws.Cells.Item(97, 14).Value = -4992.0002
$ws.Cells.Item(107, 8).Value = 5060.524
$ws.Cells.Item(107, 9).Value = 334
$ws.Cells.Item(107, 11).Value = 1002
$ws.Cells.Item(107, 13).Value = 918
$ws.Cells.Item(123, 8).Value = 2500.2068
$ws.Cells.Item(123, 9).Value = 1467.5
$ws.Cells.Item(123, 11).Value = 4402.5
$ws.Cells.Item(123, 13).Value = -1952.5
$ws.Cells.Item(124, 8).Value = 1237.1666
$ws.Cells.Item(124, 9).Value = 500
$ws.Cells.Item(124, 10).Value = 1384.6
$ws.Cells.Item(124, 11).Value = 1500
$ws.Cells.Item(124, 12).Value = 4153.799999999999
$ws.Cells.Item(124, 13).Value = 3410
$ws.Cells.Item(124, 14).Value = -13973.8
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 2334.3928
$ws.Cells.Item(122, 9).Value = 2630.2354
$ws.Cells.Item(122, 11).Value = 7890.706200000001
$ws.Cells.Item(122, 13).Value = -5440.706200000001
$ws.Cells.Item(126, 8).Value = 1840.2858
$ws.Cells.Item(126, 9).Value = 1576.4
$ws.Cells.Item(126, 10).Value = 2500
$ws.Cells.Item(126, 11).Value = 4729.200000000001
$ws.Cells.Item(126, 12).Value = 7500
$ws.Cells.Item(126, 13).Value = -2259.200000000001
$ws.Cells.Item(126, 14).Value = -12440
$ws.Cells.Item(132, 8).Value = 4028.4119
$ws.Cells.Item(132, 9).Value = 4247.1
$ws.Cells.Item(132, 10).Value = 3716
$ws.Cells.Item(132, 11).Value = 12741.3
$ws.Cells.Item(132, 12).Value = 11148
$ws.Cells.Item(132, 13).Value = -10211.3
$ws.Cells.Item(132, 14).Value = -16208
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 1222.3334
$ws.Cells.Item(61, 9).Value = 1171.5714
$ws.Cells.Item(61, 10).Value = 1400
$ws.Cells.Item(61, 11).Value = 1171.5714
$ws.Cells.Item(61, 12).Value = 1400
$ws.Cells.Item(61, 13).Value = -969.5714
$ws.Cells.Item(61, 14).Value = -1804
$ws.Cells.Item(68, 8).Value = 2198.889
$ws.Cells.Item(68, 9).Value = 2141.4285
$ws.Cells.Item(68, 11).Value = 2141.4285
$ws.Cells.Item(68, 13).Value = -1392.4285
$ws.Cells.Item(71, 8).Value = 2198.889
$ws.Cells.Item(71, 9).Value = 2141.4285
$ws.Cells.Item(71, 11).Value = 10707.1425
$ws.Cells.Item(71, 13).Value = -6963.1425
$ws.Cells.Item(113, 8).Value = 1222.3334
$ws.Cells.Item(113, 9).Value = 1171.5714
$ws.Cells.Item(113, 10).Value = 1400
$ws.Cells.Item(113, 11).Value = 1171.5714
$ws.Cells.Item(113, 12).Value = 1400
$ws.Cells.Item(113, 13).Value = 998.4286
$ws.Cells.Item(113, 14).Value = -5740
$ws.Cells.Item(114, 8).Value = 45999.332
$ws.Cells.Item(114, 10).Value = 45999.332
$ws.Cells.Item(114, 12).Value = 45999.332
$ws.Cells.Item(114, 14).Value = -54677.332
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 384.0625
$ws.Cells.Item(113, 9).Value = 269.0909
$ws.Cells.Item(113, 10).Value = 637
$ws.Cells.Item(113, 11).Value = 807.2727
$ws.Cells.Item(113, 12).Value = 1911
$ws.Cells.Item(113, 13).Value = 1362.7273
$ws.Cells.Item(113, 14).Value = -6251
